# "continental us message fix"
# Update a handful of MSRP / DPHF values on Sheet1, and replace the
# placeholder text that used to live in D34 with the real MSRP number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Base MSRP (column D) corrections -------------------------------------------------
$ws.Range("D29").Value = 53100   # 9700  - GX 460            53000 -> 53100
$ws.Range("D30").Value = 55890   # 9700PM- GX 460 Premium     55790 -> 55890
$ws.Range("D31").Value = 64365   # 9710  - GX 460 Luxury      64265 -> 64365
$ws.Range("D32").Value = 86580   # 9625  - LX 570 Two-Row     86480 -> 86580
$ws.Range("D33").Value = 91580   # 9620  - LX 570 Three-Row   91480 -> 91580

# D34 previously held a placeholder text message instead of a price.
# Replace it with the real numeric MSRP and give it the same number
# format ("#,##0") used by the rest of the Base MSRP column.
$ws.Range("D34").NumberFormat = "#,##0"
$ws.Range("D34").Value = 99310

# --- DPHF (column E) corrections -------------------------------------------------------
$ws.Range("E32").Value = 1025   # LX 570 Two-Row             1295 -> 1025
$ws.Range("E33").Value = 1025   # LX 570 Three-Row           1295 -> 1025
$ws.Range("E34").Value = 1025   # LX 570 Inspiration Series  1295 -> 1025
